# Update countries & provincias Spain
# - Swap display order of three country-name pairs (their underlying data
#   rows keep their place, but the label each row shows changes, matching
#   how the source sheet re-sorts as daily case counts shift rank).
# - Refresh the "Datos actualizados..." timestamp string.
# - Refresh case-count figures (Total/Nuevos/Activos/Recuperados/Criticos/
#   MuertesHoy/Muertes) for the rows whose countries moved or whose counts
#   were updated in this data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Country label swaps -------------------------------------------------
# Canada <-> Argentina (rows 23/24)
$ws.Cells.Item(23, 1).Value = "Argentina"
$ws.Cells.Item(24, 1).Value = "Canada"

# Georgia <-> Uruguay (rows 144/145)
$ws.Cells.Item(144, 1).Value = "Uruguay"
$ws.Cells.Item(145, 1).Value = "Georgia"

# Islas Malvinas <-> Groenlandia (rows 209/210)
$ws.Cells.Item(209, 1).Value = "Groenlandia"
$ws.Cells.Item(210, 1).Value = "Islas Malvinas"

# --- 2) Timestamp update -----------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 16 de Julio de 2020 a las 02:59"

# --- 3) Updated statistics rows ---------------------------------------------
# Map of row -> (B,C,D,E,F,G,H) new values
$rowData = @{
    4   = @(3616096, 71019, 1645715, 1830276, 0, 962, 140105)   # Estados Unidos
    22  = @(165169, 5271, 71736, 87619, 0, 189, 5814)           # Colombia
    23  = @(111146, 4236, 47298, 61798, 0, 82, 2050)            # now Argentina
    24  = @(108829, 343, 72485, 27534, 0, 12, 8810)             # now Canada
    42  = @(49243, 1147, 25417, 22844, 0, 22, 982)              # Panama
    77  = @(10428, 418, 3050, 7278, 0, 4, 100)                  # Venezuela
    95  = @(5564, 46, 2830, 2585, 0, 0, 149)                    # Mauritania
    108 = @(2831, 30, 2321, 496, 0, 0, 14)                      # Maldivas
    125 = @(1668, 17, 1200, 404, 0, 0, 64)                      # Sierra Leona
    144 = @(1009, 12, 909, 69, 0, 0, 31)                        # now Uruguay
    145 = @(1004, 5, 873, 116, 0, 0, 15)                        # now Georgia
}

foreach ($row in $rowData.Keys) {
    $values = $rowData[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 2 + $i   # column B = 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
